$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 51 (pushes existing rows 51..196 down to 52..197,
# carrying their content/formatting with them).
$ws.Rows("51").Insert()

# Populate the newly inserted row 51 with the new weekly record.
$ws.Range("A51").Value = 5
$ws.Range("B51").Value = "Macroferia Regional de Talca"
$ws.Range("C51").Value = "Maule"
$ws.Range("D51").Value = 44497
$ws.Range("E51").Value = 7
$ws.Range("F51").Value = 100114014
$ws.Range("G51").Value = "Betarraga"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 650
$ws.Range("L51").Value = 650
$ws.Range("M51").Value = 650
$ws.Range("N51").Value = "`$/paquete 5 unidades"
$ws.Range("O51").Value = "Región del Maule"
$ws.Range("P51").Value = 130
$ws.Range("Q51").Value = 5
$ws.Range("R51").Value = "Hortaliza"
